$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-5: the value changes from
# serial date 45175 (2023-09-06) to serial date 45183 (2023-09-14).
# Existing cell formatting (date number format) is preserved automatically.
$newDate = Get-Date -Year 2023 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0

$ws.Range("C2:C5").Value = $newDate
